# Fruta / hortaliza, semanal
# Insert two new weekly records (rows 447 and 448) into the Piña / Vega
# Modelo de Temuco dataset, pushing the existing rows 447-476 down to
# 449-478 (dimension grows from A1:T476 to A1:T478).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 447, shifting everything
# below (447..476) down to (449..478).
$ws.Range("A447:A448").EntireRow.Insert()

# New row 447: Primera quality, 200 boxes, $/caja 12 unidades
$ws.Range("A447").Value = 10
$ws.Range("B447").Value = "Vega Modelo de Temuco"
$ws.Range("C447").Value = "La Araucanía"
$ws.Range("D447").Value = 44714
$ws.Range("E447").Value = 9
$ws.Range("F447").Value = "Fruta"
$ws.Range("G447").Value = 100108
$ws.Range("H447").Value = "Tropicales y subtropicales"
$ws.Range("I447").Value = 100108005
$ws.Range("J447").Value = "Piña"
$ws.Range("K447").Value = "Caramelo"
$ws.Range("L447").Value = "Primera"
$ws.Range("M447").Value = 200
$ws.Range("N447").Value = 20000
$ws.Range("O447").Value = 20000
$ws.Range("P447").Value = 20000
$ws.Range("Q447").Value = "$/caja 12 unidades"
$ws.Range("R447").Value = "Ecuador"
$ws.Range("S447").Value = 1667
$ws.Range("T447").Value = 12

# New row 448: Segunda quality, 300 boxes, $/caja 7 unidades
$ws.Range("A448").Value = 10
$ws.Range("B448").Value = "Vega Modelo de Temuco"
$ws.Range("C448").Value = "La Araucanía"
$ws.Range("D448").Value = 44714
$ws.Range("E448").Value = 9
$ws.Range("F448").Value = "Fruta"
$ws.Range("G448").Value = 100108
$ws.Range("H448").Value = "Tropicales y subtropicales"
$ws.Range("I448").Value = 100108005
$ws.Range("J448").Value = "Piña"
$ws.Range("K448").Value = "Caramelo"
$ws.Range("L448").Value = "Segunda"
$ws.Range("M448").Value = 300
$ws.Range("N448").Value = 10000
$ws.Range("O448").Value = 10000
$ws.Range("P448").Value = 10000
$ws.Range("Q448").Value = "$/caja 7 unidades"
$ws.Range("R448").Value = "Ecuador"
$ws.Range("S448").Value = 1429
$ws.Range("T448").Value = 7
